# alterações a peças para imprimir
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: "encaixe superior eixo Z" -> obtido = sim, and notes cleared
$ws.Range("E4").Value = "sim"
$ws.Range("G4").Clear()

# Row 5: "base para motores inferiores Z" -> obtido = sim
$ws.Range("E5").Value = "sim"

# Row 11: "suporte para motor/veio Y" -> modelo CAD = concluído
$ws.Range("F11").Value = "concluído"

# Row 12: "suporte para veio Y" -> obtido = sim, modelo CAD = concluído, notes cleared
$ws.Range("E12").Value = "sim"
$ws.Range("F12").Value = "concluído"
$ws.Range("G12").Clear()

# Row 23: "monitor" -> modelo CAD = incompleto
$ws.Range("F23").Value = "incompleto"

# Update the visible selection/scroll position to match the author's last edit
$ws.Range("F23").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
